# Applies "Atualização de bases das ligas, do dia: 09-04-2024 às 22:40"
# - Rows 73/74 swap their match data (id, teams, odds, results)
# - Rows 112/113 swap their match data (id, teams, odds, results)
# - Rows 146, 147, 149, 150 get updated odds values

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-RowValues {
    param($Sheet, $Row, $Values)
    foreach ($col in $Values.Keys) {
        $Sheet.Range("$col$Row").Value = $Values[$col]
    }
}

# ---- Row 73 (final values, swapped in from former row 74) ----
Set-RowValues $ws 73 @{
    "B"  = 7646750
    "F"  = "Perth Glory"
    "G"  = "Wellington Phoenix"
    "H"  = 3
    "I"  = 4
    "J"  = "A"
    "K"  = 2.45
    "L"  = 3.75
    "M"  = 2.55
    "N"  = 3.1
    "O"  = 3.8
    "P"  = 2.05
    "Q"  = 0.25
    "R"  = 2
    "S"  = 1.85
    "T"  = 3
    "U"  = 1.925
    "V"  = 1.925
    "W"  = -1
    "X"  = -1
    "Y"  = 1.05
    "Z"  = -1
    "AA" = 0.8500000000000001
    "AB" = 0.925
    "AC" = -1
}

# ---- Row 74 (final values, swapped in from former row 73) ----
Set-RowValues $ws 74 @{
    "B"  = 7646749
    "F"  = "Brisbane Roar"
    "G"  = "Newcastle Jets"
    "H"  = 3
    "I"  = 2
    "J"  = "H"
    "K"  = 1.909
    "L"  = 4
    "M"  = 3.4
    "N"  = 2.4
    "O"  = 4
    "P"  = 2.6
    "Q"  = 0
    "R"  = 1.83
    "S"  = 2.07
    "T"  = 3.25
    "U"  = 1.9
    "V"  = 1.95
    "W"  = 1.4
    "X"  = -1
    "Y"  = -1
    "Z"  = 0.8300000000000001
    "AA" = -1
    "AB" = 0.8999999999999999
    "AC" = -1
}

# ---- Row 112 (final values, swapped in from former row 113) ----
Set-RowValues $ws 112 @{
    "B"  = 7127379
    "F"  = "Melbourne Victory"
    "G"  = "Central Coast Mariners"
    "H"  = 0
    "I"  = 1
    "J"  = "A"
    "K"  = 1.95
    "L"  = 3.6
    "M"  = 3.8
    "N"  = 1.909
    "O"  = 3.6
    "P"  = 4
    "Q"  = -0.5
    "R"  = 1.9
    "S"  = 1.95
    "T"  = 2.75
    "U"  = 1.925
    "V"  = 1.925
    "W"  = -1
    "X"  = -1
    "Y"  = 3
    "Z"  = -1
    "AA" = 0.95
    "AB" = -1
    "AC" = 0.925
}

# ---- Row 113 (final values, swapped in from former row 112) ----
Set-RowValues $ws 113 @{
    "B"  = 7127376
    "F"  = "Newcastle Jets"
    "G"  = "Macarthur FC"
    "H"  = 2
    "I"  = 2
    "J"  = "D"
    "K"  = 1.95
    "L"  = 4
    "M"  = 3.4
    "N"  = 1.909
    "O"  = 4.2
    "P"  = 3.6
    "Q"  = -0.5
    "R"  = 1.89
    "S"  = 2.01
    "T"  = 3.5
    "U"  = 1.95
    "V"  = 1.9
    "W"  = -1
    "X"  = 3.2
    "Y"  = -1
    "Z"  = -1
    "AA" = 1.01
    "AB" = 0.95
    "AC" = -1
}

# ---- Row 146: odds updates ----
Set-RowValues $ws 146 @{
    "N" = 3.6
    "Q" = 0.5
    "R" = 1.83
    "S" = 2.07
}

# ---- Row 147: odds updates ----
Set-RowValues $ws 147 @{
    "O" = 4.333
    "P" = 3.8
    "R" = 2.05
    "S" = 1.85
    "U" = 1.85
    "V" = 2
}

# ---- Row 149: odds updates ----
Set-RowValues $ws 149 @{
    "R" = 2.01
    "S" = 1.89
}

# ---- Row 150: odds updates ----
Set-RowValues $ws 150 @{
    "Q" = -1
    "R" = 2.07
    "S" = 1.83
}
